# docs: added train_2022 components
#
# Duplicate the "train" sheet into a new "train_2022" sheet (placed right
# after "train"), tweak the new sheet's first phase duration (B2: 180 -> 120),
# make it the active/selected sheet, and hide the older per-session detail
# sheets that are no longer the primary view.

$wb = $excel.ActiveWorkbook

# 1) Duplicate "train" -> new sheet placed immediately after it.
$trainSheet = $wb.Worksheets.Item("train")
$trainSheet.Copy($null, $trainSheet)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "train_2022"

# 2) Update the first baseline phase duration for the 2022 protocol.
$newSheet.Range("B2").Value = 120

# 3) Make the new sheet the active tab, with the same kind of selection
#    state Excel leaves behind after editing near the bottom of the table.
$newSheet.Select()
$newSheet.Range("B42").Select()

# 4) Hide the older/legacy sheets now that train_2022 is the primary sheet.
$wb.Worksheets.Item("tone").Visible = $false
$wb.Worksheets.Item("extinction").Visible = $false
$wb.Worksheets.Item("Pav_app").Visible = $false
$wb.Worksheets.Item("cs_response").Visible = $false
$wb.Worksheets.Item("shock_response").Visible = $false

$lastSheet = $wb.Worksheets.Item("cs_response_2")
$lastSheet.Activate()
$lastSheet.Range("A1:D1").Select()
$lastSheet.Visible = $false

# 5) Leave focus back on the newly-added sheet, matching the saved file's
#    active tab.
$newSheet.Activate()
